$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.352
$ws.Range("B4").Value = 5.125
$ws.Range("A6").Value = -22.468
$ws.Range("D6").Value = -8.484999999999999
$ws.Range("A7").Value = -20.978
$ws.Range("D7").Value = -8.128
$ws.Range("A8").Value = -21.675
$ws.Range("B8").Value = 6.418000000000001
$ws.Range("D8").Value = -8.351000000000001
$ws.Range("B9").Value = 6.456
$ws.Range("D10").Value = -8.071000000000002
$ws.Range("B12").Value = 5.995000000000001
$ws.Range("D13").Value = -8.181999999999999
$ws.Range("D14").Value = -8.059000000000001
$ws.Range("A16").Value = -21.254
$ws.Range("D16").Value = -8.381
$ws.Range("B17").Value = 5.39
$ws.Range("B18").Value = 5.789999999999999
$ws.Range("B19").Value = 7.311
$ws.Range("A20").Value = -21.902
$ws.Range("B20").Value = 5.256
$ws.Range("A21").Value = -20.186
$ws.Range("B26").Value = 5.972
$ws.Range("A28").Value = -21.604
$ws.Range("A29").Value = -21.5
$ws.Range("A30").Value = -21.275
$ws.Range("D30").Value = -7.496
$ws.Range("B31").Value = 6.718999999999999
$ws.Range("A32").Value = -21.439
$ws.Range("D37").Value = -8.244999999999999
$ws.Range("B39").Value = 7.092000000000001
$ws.Range("A40").Value = -21.203
$ws.Range("B40").Value = 7.042
$ws.Range("D40").Value = -8.301
$ws.Range("B41").Value = 6.462000000000001
$ws.Range("B42").Value = 6.382000000000001
$ws.Range("B43").Value = 6.786
$ws.Range("D44").Value = -7.371
$ws.Range("A46").Value = -21.489
$ws.Range("B47").Value = 6.051
$ws.Range("B48").Value = 6.194
$ws.Range("A51").Value = -21.176
$ws.Range("A52").Value = -21.646
$ws.Range("B54").Value = 5.278
$ws.Range("A57").Value = -21.806
$ws.Range("A59").Value = -22.257
$ws.Range("A62").Value = -21.897
$ws.Range("B62").Value = 5.258000000000001
$ws.Range("B63").Value = 5.484
$ws.Range("B64").Value = 6.006
$ws.Range("A66").Value = -21.482
$ws.Range("D70").Value = -6.98
$ws.Range("A73").Value = -21.367
$ws.Range("A74").Value = -20.673
$ws.Range("B76").Value = 6.255999999999998
$ws.Range("A77").Value = -21.503
$ws.Range("B81").Value = 5.286
$ws.Range("B84").Value = 5.924000000000001
$ws.Range("B89").Value = 6.023000000000001
$ws.Range("D89").Value = -8.347000000000001
$ws.Range("D91").Value = -7.613
$ws.Range("A92").Value = -21.566
$ws.Range("D93").Value = -6.915000000000001
$ws.Range("B94").Value = 6.503
$ws.Range("D98").Value = -7.187
$ws.Range("A100").Value = -22.031
